# Auto-generated edit script applying numeric updates to Alexander_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1772.1666
$ws.Range("I8").Value = 78.333336
$ws.Range("J8").Value = 3466
$ws.Range("K8").Value = 235.000008
$ws.Range("L8").Value = 10398
$ws.Range("M8").Value = -96.00000800000001
$ws.Range("N8").Value = -10676

$ws.Range("H43").Value = 824.44446
$ws.Range("I43").Value = 455
$ws.Range("J43").Value = 1120
$ws.Range("K43").Value = 455
$ws.Range("L43").Value = 1120
$ws.Range("M43").Value = -386
$ws.Range("N43").Value = -1258

$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 6000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -6540

$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 6000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -7872

$ws.Range("H97").Value = 1299.5416
$ws.Range("J97").Value = 1258.6364
$ws.Range("L97").Value = 3775.9092
$ws.Range("N97").Value = -4767.9092

$ws.Range("H129").Value = 1130.7826
$ws.Range("I129").Value = 335.2
$ws.Range("J129").Value = 2622.5
$ws.Range("K129").Value = 1005.6
$ws.Range("L129").Value = 7867.5
$ws.Range("M129").Value = 3994.4
$ws.Range("N129").Value = -17867.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19121.955
$ws.Range("I32").Value = 3937.6667
$ws.Range("K32").Value = 3937.6667
$ws.Range("M32").Value = -3650.6667

$ws.Range("H110").Value = 2678.077
$ws.Range("I110").Value = 3129.85
$ws.Range("J110").Value = 1172.1666
$ws.Range("K110").Value = 3129.85
$ws.Range("L110").Value = 1172.1666
$ws.Range("M110").Value = -1084.85
$ws.Range("N110").Value = -5262.1666

$ws.Range("H122").Value = 2168.4736
$ws.Range("I122").Value = 1900.3077
$ws.Range("J122").Value = 2749.5
$ws.Range("K122").Value = 5700.9231
$ws.Range("L122").Value = 8248.5
$ws.Range("M122").Value = -3250.9231
$ws.Range("N122").Value = -13148.5

$ws.Range("H132").Value = 2590.5806
$ws.Range("I132").Value = 2859.15
$ws.Range("J132").Value = 2102.2727
$ws.Range("K132").Value = 8577.450000000001
$ws.Range("L132").Value = 6306.8181
$ws.Range("M132").Value = -6047.450000000001
$ws.Range("N132").Value = -11366.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3246.973
$ws.Range("I107").Value = 3425.4644
$ws.Range("J107").Value = 2691.6667
$ws.Range("K107").Value = 3425.4644
$ws.Range("L107").Value = 2691.6667
$ws.Range("M107").Value = -1505.4644
$ws.Range("N107").Value = -6531.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1846.8889
$ws.Range("I31").Value = 1251.5518
$ws.Range("J31").Value = 4313.2856
$ws.Range("K31").Value = 1251.5518
$ws.Range("L31").Value = 4313.2856
$ws.Range("M31").Value = -956.5518
$ws.Range("N31").Value = -4903.2856

$ws.Range("H34").Value = 1846.8889
$ws.Range("I34").Value = 1251.5518
$ws.Range("J34").Value = 4313.2856
$ws.Range("K34").Value = 1251.5518
$ws.Range("L34").Value = 4313.2856
$ws.Range("M34").Value = -1049.5518
$ws.Range("N34").Value = -4717.2856

$ws.Range("H105").Value = 5823.4165
$ws.Range("I105").Value = 5914.4287
$ws.Range("J105").Value = 5696
$ws.Range("K105").Value = 5914.4287
$ws.Range("L105").Value = 5696
$ws.Range("M105").Value = -4167.4287
$ws.Range("N105").Value = -9190

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 805.8182
$ws.Range("J5").Value = 1150
$ws.Range("L5").Value = 3450
$ws.Range("N5").Value = -3674

$ws.Range("H8").Value = 83.15385000000001
$ws.Range("I8").Value = 83.15385000000001
$ws.Range("K8").Value = 249.46155
$ws.Range("M8").Value = -110.46155

$ws.Range("H81").Value = 90910940
$ws.Range("I81").Value = 540
$ws.Range("J81").Value = 111113256
$ws.Range("K81").Value = 1620
$ws.Range("L81").Value = 333339768
$ws.Range("M81").Value = -497
$ws.Range("N81").Value = -333342014

$ws.Range("H84").Value = 90910940
$ws.Range("I84").Value = 540
$ws.Range("J84").Value = 111113256
$ws.Range("K84").Value = 4860
$ws.Range("L84").Value = 1000019304
$ws.Range("M84").Value = 756
$ws.Range("N84").Value = -1000030536

$ws.Range("H109").Value = 9397.706
$ws.Range("I109").Value = 11836.1
$ws.Range("K109").Value = 35508.3
$ws.Range("M109").Value = -34468.3

$ws.Range("H129").Value = 169191.33
$ws.Range("I129").Value = 334953.34
$ws.Range("J129").Value = 3429.3333
$ws.Range("K129").Value = 1004860.02
$ws.Range("L129").Value = 10287.9999
$ws.Range("M129").Value = -999860.02
$ws.Range("N129").Value = -20287.9999

$ws.Range("H135").Value = 805.8182
$ws.Range("J135").Value = 1150
$ws.Range("L135").Value = 10350
$ws.Range("N135").Value = -15420

$ws.Range("H140").Value = 32969.5
$ws.Range("I140").Value = 63439.75
$ws.Range("J140").Value = 2499.25
$ws.Range("K140").Value = 190319.25
$ws.Range("L140").Value = 7497.75
$ws.Range("M140").Value = -185139.25
$ws.Range("N140").Value = -17857.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 38440
$ws.Range("J123").Value = 38440
$ws.Range("L123").Value = 38440
$ws.Range("N123").Value = -43340

$ws.Range("H132").Value = 2667.158
$ws.Range("I132").Value = 2175
$ws.Range("J132").Value = 3510.8572
$ws.Range("K132").Value = 6525
$ws.Range("L132").Value = 10532.5716
$ws.Range("M132").Value = -3995
$ws.Range("N132").Value = -15592.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 733.3333
$ws.Range("I6").Value = 200
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 200
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = -85
$ws.Range("N6").Value = -1230

$ws.Range("H41").Value = 5111.4
$ws.Range("J41").Value = 5111.4
$ws.Range("L41").Value = 5111.4
$ws.Range("N41").Value = -5891.4

$ws.Range("H64").Value = 23664.4
$ws.Range("J64").Value = 23664.4
$ws.Range("L64").Value = 23664.4
$ws.Range("N64").Value = -24160.4

$ws.Range("H67").Value = 23664.4
$ws.Range("J67").Value = 23664.4
$ws.Range("L67").Value = 23664.4
$ws.Range("N67").Value = -25380.4

$ws.Range("H136").Value = 1129.5952
$ws.Range("I136").Value = 899.25
$ws.Range("J136").Value = 2511.6667
$ws.Range("K136").Value = 2697.75
$ws.Range("L136").Value = 7535.000100000001
$ws.Range("M136").Value = -147.75
$ws.Range("N136").Value = -12635.0001
